# Refresh cryptos.xlsx Price (D) and Volume(1h) (E) columns with latest
# coinranking.com snapshot values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.172.13"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.677.06"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.50"
$ws.Range("E5").Value = "  -3.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5286"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2680"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06289"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.31"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07511"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "1.680.13"
$ws.Range("E12").Value = "  +7.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.483"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5656"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008113"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.19"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "26.219.40"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.850"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.76"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.201"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.16"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1261"
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.607"
$ws.Range("E26").Value = "  -3.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.86"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06474"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.344"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.524"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.482"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.648"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6086"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "1.100.78"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8664"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.006"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.95"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "1.830.18"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.79"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05268"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.994"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4269"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.951"
$ws.Range("E51").Value = "  -2.11%  "
